# Scheduled data refresh: update currentAveragePrice* / LevePrice* / LeveProfit*
# columns (H-N) across all 8 Leve sheets with latest Market Board pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$updates = @(
    @("H15", 882.43286),
    @("I15", 882.43286),
    @("K15", 2647.29858),
    @("M15", -2478.29858),
    @("H48", 4575),
    @("J48", 4575),
    @("L48", 13725),
    @("N48", -14309),
    @("H56", 4575),
    @("J56", 4575),
    @("L56", 13725),
    @("N56", -14793),
    @("H62", 31251956),
    @("I62", 38463316),
    @("K62", 38463316),
    @("M62", -38462692),
    @("H65", 31251956),
    @("I65", 38463316),
    @("K65", 192316580),
    @("M65", -192313460),
    @("H69", 23258350),
    @("I69", 5500),
    @("J69", 24392636),
    @("K69", 16500),
    @("L69", 73177908),
    @("M69", -15626),
    @("N69", -73179656),
    @("H72", 23258350),
    @("I72", 5500),
    @("J72", 24392636),
    @("K72", 49500),
    @("L72", 219533724),
    @("M72", -45132),
    @("N72", -219542460),
    @("H96", 783.1),
    @("I96", 667.4286),
    @("J96", 1053),
    @("K96", 2002.2858),
    @("L96", 3159),
    @("M96", -629.2857999999999),
    @("N96", -5905),
    @("H98", 1261),
    @("I98", 1261),
    @("K98", 1261),
    @("M98", 237),
    @("H100", 1600),
    @("I100", 0),
    @("J100", 1600),
    @("K100", 0),
    @("L100", 1600),
    @("M100", $null),
    @("N100", -2682),
    @("H101", 3290.6),
    @("I101", 634),
    @("J101", 5061.6665),
    @("K101", 1902),
    @("L101", 15184.9995),
    @("M101", -280),
    @("N101", -18428.9995),
    @("H113", 1988.1052),
    @("I113", 1777),
    @("J113", 1999.8334),
    @("K113", 1777),
    @("L113", 1999.8334),
    @("M113", 1477),
    @("N113", -8507.8334),
    @("H122", 1261),
    @("I122", 1261),
    @("K122", 3783),
    @("M122", -1333),
    @("H125", 2005.8),
    @("J125", 2018),
    @("L125", 18162),
    @("N125", -23082),
    @("H127", 642.8),
    @("I127", 499.25),
    @("J127", 1217),
    @("K127", 1497.75),
    @("L127", 3651),
    @("M127", 3462.25),
    @("N127", -13571),
    @("H138", 3205.4211),
    @("I138", 726.8421),
    @("J138", 5684),
    @("K138", 2180.5263),
    @("L138", 17052),
    @("M138", 2959.4737),
    @("N138", -27332)
)
foreach ($u in $updates) {
    if ($null -eq $u[1]) {
        $ws.Range($u[0]).ClearContents()
    } else {
        $ws.Range($u[0]).Value = $u[1]
    }
}

$ws = $wb.Sheets.Item("ARM")
$updates = @(
    @("H2", 1606.04),
    @("I2", 1684.1177),
    @("J2", 1440.125),
    @("K2", 1684.1177),
    @("L2", 1440.125),
    @("M2", -1571.1177),
    @("N2", -1666.125),
    @("H6", 7847.8184),
    @("J6", 11428.143),
    @("L6", 11428.143),
    @("N6", -11774.143),
    @("H32", 4259.5244),
    @("I32", 1888.5193),
    @("K32", 1888.5193),
    @("M32", -1601.5193),
    @("H45", 3113.5),
    @("I45", 3084.75),
    @("J45", 3156.625),
    @("K45", 3084.75),
    @("L45", 3156.625),
    @("M45", -2707.75),
    @("N45", -3910.625),
    @("H116", 1606.04),
    @("I116", 1684.1177),
    @("J116", 1440.125),
    @("K116", 1684.1177),
    @("L116", 1440.125),
    @("M116", 609.8823),
    @("N116", -6028.125),
    @("H122", 1373212.6),
    @("I122", 1425913.1),
    @("K122", 4277739.300000001),
    @("M122", -4275289.300000001)
)
foreach ($u in $updates) {
    if ($null -eq $u[1]) {
        $ws.Range($u[0]).ClearContents()
    } else {
        $ws.Range($u[0]).Value = $u[1]
    }
}

$ws = $wb.Sheets.Item("BSM")
$updates = @(
    @("H3", 1606.04),
    @("I3", 1684.1177),
    @("J3", 1440.125),
    @("K3", 1684.1177),
    @("L3", 1440.125),
    @("M3", -1570.1177),
    @("N3", -1668.125),
    @("H125", 45780),
    @("J125", 45780),
    @("L125", 45780),
    @("N125", -55620)
)
foreach ($u in $updates) {
    if ($null -eq $u[1]) {
        $ws.Range($u[0]).ClearContents()
    } else {
        $ws.Range($u[0]).Value = $u[1]
    }
}

$ws = $wb.Sheets.Item("CRP")
$updates = @(
    @("H2", 10142.143),
    @("I2", 5498.75),
    @("J2", 16333.333),
    @("K2", 5498.75),
    @("L2", 16333.333),
    @("M2", -5385.75),
    @("N2", -16559.333),
    @("H3", 28000),
    @("J3", 28000),
    @("L3", 28000),
    @("N3", -28226),
    @("H4", 10000),
    @("I4", 0),
    @("J4", 10000),
    @("K4", 0),
    @("L4", 10000),
    @("M4", $null),
    @("N4", -10224),
    @("H6", 7201200),
    @("I6", 12000000),
    @("K6", 12000000),
    @("M6", -11999887),
    @("H7", 124.1875),
    @("I7", 45),
    @("K7", 45),
    @("M7", 68),
    @("H10", 3935),
    @("I10", 513.3333),
    @("J10", 14200),
    @("K10", 513.3333),
    @("L10", 14200),
    @("M10", -374.3333),
    @("N10", -14478),
    @("H11", 20000),
    @("I11", 0),
    @("J11", 20000),
    @("K11", 0),
    @("L11", 20000),
    @("M11", $null),
    @("N11", -20280),
    @("H12", 7234.1665),
    @("I12", 851.25),
    @("J12", 20000),
    @("K12", 851.25),
    @("L12", 20000),
    @("M12", -681.25),
    @("N12", -20340),
    @("H13", 28000),
    @("I13", 0),
    @("J13", 28000),
    @("K13", 0),
    @("L13", 28000),
    @("M13", $null),
    @("N13", -28278)
)
foreach ($u in $updates) {
    if ($null -eq $u[1]) {
        $ws.Range($u[0]).ClearContents()
    } else {
        $ws.Range($u[0]).Value = $u[1]
    }
}

$ws = $wb.Sheets.Item("CUL")
$updates = @(
    @("H80", 1870.5714),
    @("I80", 1850.5),
    @("J80", 1878.6),
    @("K80", 5551.5),
    @("L80", 5635.799999999999),
    @("M80", -4615.5),
    @("N80", -7507.799999999999),
    @("H83", 1870.5714),
    @("I83", 1850.5),
    @("J83", 1878.6),
    @("K83", 16654.5),
    @("L83", 16907.4),
    @("M83", -11974.5),
    @("N83", -26267.4)
)
foreach ($u in $updates) {
    if ($null -eq $u[1]) {
        $ws.Range($u[0]).ClearContents()
    } else {
        $ws.Range($u[0]).Value = $u[1]
    }
}

$ws = $wb.Sheets.Item("GSM")
$updates = @(
    @("H102", 2048.5417),
    @("I102", 2104.0625),
    @("J102", 1937.5),
    @("K102", 2104.0625),
    @("L102", 1937.5),
    @("M102", -482.0625),
    @("N102", -5181.5),
    @("H107", 299.26315),
    @("I107", 234.5),
    @("J107", 480.6),
    @("K107", 234.5),
    @("L107", 480.6),
    @("M107", 1685.5),
    @("N107", -4320.6),
    @("H122", 9092846),
    @("I122", 12501850),
    @("J122", 2169.3333),
    @("K122", 37505550),
    @("L122", 6507.999899999999),
    @("M122", -37503100),
    @("N122", -11407.9999),
    @("H132", 3028.8462),
    @("I132", 2426.2354),
    @("K132", 7278.706200000001),
    @("M132", -4748.706200000001)
)
foreach ($u in $updates) {
    if ($null -eq $u[1]) {
        $ws.Range($u[0]).ClearContents()
    } else {
        $ws.Range($u[0]).Value = $u[1]
    }
}

$ws = $wb.Sheets.Item("LTW")
$updates = @(
    @("H40", 6554.6665),
    @("I40", 8729.200000000001),
    @("J40", 5001.4287),
    @("K40", 8729.200000000001),
    @("L40", 5001.4287),
    @("M40", -8593.200000000001),
    @("N40", -5273.4287),
    @("H93", 25942.25),
    @("I93", 33921.668),
    @("J93", 2004),
    @("K93", 33921.668),
    @("L93", 2004),
    @("M93", -32673.668),
    @("N93", -4500),
    @("H136", 2031.8),
    @("I136", 1098.4615),
    @("J136", 3042.9167),
    @("K136", 3295.3845),
    @("L136", 9128.750100000001),
    @("M136", -745.3844999999997),
    @("N136", -14228.7501)
)
foreach ($u in $updates) {
    if ($null -eq $u[1]) {
        $ws.Range($u[0]).ClearContents()
    } else {
        $ws.Range($u[0]).Value = $u[1]
    }
}

$ws = $wb.Sheets.Item("WVR")
$updates = @(
    @("H62", 33334834),
    @("I62", 50001250),
    @("J62", 2000),
    @("K62", 50001250),
    @("L62", 2000),
    @("M62", -50000626),
    @("N62", -3248),
    @("H65", 33334834),
    @("I65", 50001250),
    @("J65", 2000),
    @("K65", 250006250),
    @("L65", 10000),
    @("M65", -250003130),
    @("N65", -16240),
    @("H96", 2419.923),
    @("I96", 1824.1666),
    @("J96", 2930.5715),
    @("K96", 1824.1666),
    @("L96", 2930.5715),
    @("M96", -451.1666),
    @("N96", -5676.5715),
    @("H122", 2789.1667),
    @("I122", 2317.5),
    @("J122", 3126.0715),
    @("K122", 6952.5),
    @("L122", 9378.2145),
    @("M122", -4502.5),
    @("N122", -14278.2145),
    @("H126", 1040.3077),
    @("I126", 1040.3077),
    @("K126", 3120.9231),
    @("M126", -650.9231)
)
foreach ($u in $updates) {
    if ($null -eq $u[1]) {
        $ws.Range($u[0]).ClearContents()
    } else {
        $ws.Range($u[0]).Value = $u[1]
    }
}

